# vc-5 feature/add creadit and whitelisting
# Replace the msisdn list in column A with the new set of numbers (rows 2-21),
# extending the sheet from A1:A8 to A1:A21. All values are stored as text
# (the column already used numberStoredAsText-ignored text values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "7014118238",
    "7075812222",
    "9885861677",
    "9902459657",
    "8970511445",
    "9535114669",
    "9901396041",
    "9845704305",
    "9164558164",
    "9844329150",
    "6363179872",
    "9916725929",
    "9481270887",
    "7760493586",
    "9448107102",
    "9844060947",
    "9663093906",
    "9632535125",
    "9480230542",
    "9008150443"
)

# Make sure the whole target range is formatted as Text *before* writing,
# so the numeric-looking strings are not reinterpreted as numbers.
$ws.Range("A1:A21").NumberFormat = "@"

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

Write-Host "Updated A2:A21 with new msisdn values"
